# تعديل تلقائي في شيت Card5 by admin at 2025-11-23 09:15:54
#
# Card5 worksheet: tidy up the "Serviced by" header (drop the trailing
# space) and populate the "O" column's empty placeholder cells with the
# same literal "nan" text used throughout the rest of the table. Also
# records a service note ("تم سن الفلاتس") for row 8, where the flats
# servicing took place on 23/11/2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card5")

# O1 header: "Serviced by " -> "Serviced by" (trailing space removed)
$ws.Range("O1").Value = "Serviced by"

# O column placeholders that were blank become the literal text "nan",
# matching the rest of the sheet's empty-value convention.
$ws.Range("O2").Value = "nan"
$ws.Range("O3").Value = "nan"
$ws.Range("O4").Value = "nan"
$ws.Range("O5").Value = "nan"
$ws.Range("O6").Value = "nan"
$ws.Range("O7").Value = "nan"

# Row 8: log the flats-sharpening service event in the Event column (M),
# keep the Correction column (N) as "nan", and fill O8 with "nan" too.
$ws.Range("M8").Value = "تم سن الفلاتس"
$ws.Range("N8").Value = "nan"
$ws.Range("O8").Value = "nan"

$ws.Range("O9").Value = "nan"
$ws.Range("O10").Value = "nan"
$ws.Range("O11").Value = "nan"
$ws.Range("O12").Value = "nan"
$ws.Range("O13").Value = "nan"
